$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = 0.5566666666666666
$ws.Range("J2").Value = 0.01333333333333333
$ws.Range("P2").Value = 0.14
$ws.Range("S2").Value = 0.09
$ws.Range("C3").Value = 0.02298850574712644
$ws.Range("J3").Value = 0.06896551724137931
$ws.Range("P3").Value = 0.7068965517241379
$ws.Range("S3").Value = 0.2011494252873563
$ws.Range("B6").Value = 0.08658008658008658
$ws.Range("D6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.1168831168831169
$ws.Range("J6").Value = 0.2424242424242424
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.2034632034632035
$ws.Range("R6").Value = 0.05627705627705628
$ws.Range("S6").Value = 0.2813852813852814
$ws.Range("B7").Value = 0.09142857142857143
$ws.Range("D7").Value = 0.02857142857142857
$ws.Range("F7").Value = 0.05142857142857143
$ws.Range("J7").Value = 0.1371428571428571
$ws.Range("O7").Value = 0.04
$ws.Range("Q7").Value = 0.2171428571428571
$ws.Range("R7").Value = 0.1028571428571429
$ws.Range("S7").Value = 0.3314285714285714
$ws.Range("B8").Value = 0.1271551724137931
$ws.Range("D8").Value = 0.01939655172413793
$ws.Range("F8").Value = 0.06465517241379311
$ws.Range("J8").Value = 0.09267241379310345
$ws.Range("O8").Value = 0.02155172413793104
$ws.Range("Q8").Value = 0.1767241379310345
$ws.Range("R8").Value = 0.09698275862068965
$ws.Range("S8").Value = 0.4008620689655172
$ws.Range("B9").Value = 0.09714285714285714
$ws.Range("D9").Value = 0.01714285714285714
$ws.Range("F9").Value = 0.06857142857142857
$ws.Range("J9").Value = 0.1657142857142857
$ws.Range("O9").Value = 0.01714285714285714
$ws.Range("Q9").Value = 0.1828571428571429
$ws.Range("R9").Value = 0.09714285714285714
$ws.Range("S9").Value = 0.3542857142857143
$ws.Range("B10").Value = 0.1072013093289689
$ws.Range("D10").Value = 0.0204582651391162
$ws.Range("F10").Value = 0.0662847790507365
$ws.Range("J10").Value = 0.1039279869067103
$ws.Range("O10").Value = 0.0237315875613748
$ws.Range("Q10").Value = 0.2315875613747954
$ws.Range("R10").Value = 0.1014729950900164
$ws.Range("S10").Value = 0.3453355155482815
$ws.Range("F11").Value = 0.004166666666666667
$ws.Range("G11").Value = 0.1458333333333333
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.5541666666666667
$ws.Range("S11").Value = 0.008333333333333333
$ws.Range("G12").Value = 0.7971014492753623
$ws.Range("J12").Value = 0.1449275362318841
$ws.Range("K12").Value = 0.007246376811594203
$ws.Range("L12").Value = 0.02173913043478261
$ws.Range("S12").Value = 0.02898550724637681
$ws.Range("G13").Value = 0.7708333333333334
$ws.Range("J13").Value = 0.2291666666666667
$ws.Range("F15").Value = 0.01244813278008299
$ws.Range("H15").Value = 0.1327800829875519
$ws.Range("I15").Value = 0.05809128630705394
$ws.Range("J15").Value = 0.3526970954356847
$ws.Range("K15").Value = 0.0912863070539419
$ws.Range("O15").Value = 0.1037344398340249
$ws.Range("S15").Value = 0.2489626556016598
$ws.Range("F16").Value = 0.02688172043010753
$ws.Range("H16").Value = 0.1774193548387097
$ws.Range("I16").Value = 0.05913978494623656
$ws.Range("J16").Value = 0.4408602150537634
$ws.Range("K16").Value = 0.08064516129032258
$ws.Range("M16").Value = 0.01075268817204301
$ws.Range("O16").Value = 0.06989247311827956
$ws.Range("S16").Value = 0.1344086021505376
$ws.Range("F17").Value = 0.01663201663201663
$ws.Range("H17").Value = 0.20997920997921
$ws.Range("I17").Value = 0.07276507276507277
$ws.Range("J17").Value = 0.4261954261954262
$ws.Range("K17").Value = 0.06237006237006237
$ws.Range("M17").Value = 0.03742203742203742
$ws.Range("N17").Value = 0.002079002079002079
$ws.Range("O17").Value = 0.05197505197505198
$ws.Range("S17").Value = 0.1205821205821206
$ws.Range("F18").Value = 0.01376146788990826
$ws.Range("H18").Value = 0.1972477064220184
$ws.Range("I18").Value = 0.09174311926605505
$ws.Range("J18").Value = 0.4311926605504587
$ws.Range("K18").Value = 0.07798165137614679
$ws.Range("M18").Value = 0.01834862385321101
$ws.Range("N18").Value = 0.004587155963302753
$ws.Range("O18").Value = 0.06422018348623854
$ws.Range("S18").Value = 0.1009174311926606
$ws.Range("F19").Value = 0.02272727272727273
$ws.Range("H19").Value = 0.2298951048951049
$ws.Range("I19").Value = 0.08566433566433566
$ws.Range("J19").Value = 0.368006993006993
$ws.Range("K19").Value = 0.09178321678321678
$ws.Range("M19").Value = 0.02185314685314685
$ws.Range("N19").Value = 0.001748251748251748
$ws.Range("O19").Value = 0.07517482517482517
$ws.Range("S19").Value = 0.1031468531468532